# Update column F ("dSF") values on Sheet1 per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = -3
    3  = -2
    4  = -3
    6  = 1
    7  = 5
    8  = -6
    9  = -7
    10 = -1
    11 = -4
    12 = -2
    13 = 6
    14 = -1
    15 = 1
    16 = 2
    18 = -4
    19 = -1
    20 = -3
    21 = 2
    22 = 3
    24 = 1
    26 = -2
    27 = -1
    31 = -1
    32 = -2
    36 = 6
    37 = -3
    38 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
